$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlLeft = -4131
$xlTop = -4160
$xlBottom = -4107
$xlPasteFormats = -4122

# --- Row 2: Tuan (week number) ---
$ws.Range("B2").Value = 8
$ws.Range("B2:H6").VerticalAlignment = $xlBottom

# --- Row 3: Ngay bat dau (start date) ---
$ws.Range("B3").NumberFormat = "mm-dd-yy"
$ws.Range("B3").Value = "11/24/2025"

# --- Row 4: Ngay ket thuc (end date) ---
# copy the date style from B3 so both rows share a single new style entry
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial($xlPasteFormats)
$ws.Range("B4").Value = "11/23/2025"

# --- Row 5: Ho ten (full name) ---
$ws.Range("B5").Value = "Lại Việt Anh"

# --- Row 6: MSSV (student id) ---
$ws.Range("B6").Value = 179066

# --- Row 7: De tai (topic) - now wraps, left/top aligned instead of centered ---
$ws.Range("B7").Value = "Xây dựng mô hình giám sát và điều khiển nhà thông minh"
$ws.Range("B7:H7").HorizontalAlignment = $xlLeft
$ws.Range("B7:H7").VerticalAlignment = $xlTop

# --- Row 8: Cong viec da hoan thanh (completed work) ---
$ws.Range("B8").Value = "thiết kế hệ thống"
$ws.Range("B8:H8").VerticalAlignment = $xlBottom
$ws.Rows(8).RowHeight = 15.6

# --- Row 9: Cong viec du kien tuan toi (planned work) ---
$ws.Range("B9").Value = "thiết kế hệ thống"
$ws.Range("B9:H9").VerticalAlignment = $xlBottom
$ws.Rows(9).RowHeight = 15.6

# --- Selection ---
$ws.Range("B4:H4").Select | Out-Null
